# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" worksheets to reflect newly generated data
# (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13899
$ws1.Range("F6").Value = 497
$ws1.Range("F7").Value = 1213
$ws1.Range("F8").Value = 1038
$ws1.Range("F9").Value = 13906
$ws1.Range("F10").Value = 14847
$ws1.Range("F12").Value = 6
$ws1.Range("F20").Value = 23
$ws1.Range("F26").Value = 5768
$ws1.Range("F29").Value = 5427
$ws1.Range("F30").Value = 52
$ws1.Range("F32").Value = 291

# --- Sheet: 全部类型 -----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13899
$ws4.Range("F7").Value = 497
$ws4.Range("F8").Value = 1213
$ws4.Range("F9").Value = 1038
$ws4.Range("F10").Value = 13906
$ws4.Range("F11").Value = 14847
$ws4.Range("F13").Value = 6
$ws4.Range("F21").Value = 23
$ws4.Range("F27").Value = 5768
$ws4.Range("F30").Value = 5427
$ws4.Range("F31").Value = 52
$ws4.Range("F33").Value = 291
